$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, corrected "StatQuery" (column C) text shared by the Cases/Samples/Files rows.
# The old broken "all_studies" query is replaced with a query that counts
# Programs/Studies/Cases/Samples/Case Files/Study Files for the breed filter.
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['French Bulldog']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Update the view: zoomed to 85%, scrolled/selected on row 4.
$ws.Activate()
$ws.Range("B4").Select()
$excel.ActiveWindow.Zoom = 85
